# Add a new "2023" column (L) to the corruption-index data table, mirroring
# the existing "2022" column (K) for layout/formatting, and bump a few row
# heights that Excel re-flowed once the extra column was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for column L (2023) -----------------------------------------
$ws.Range("L4").Value  = 2023
$ws.Range("L5").Value  = 22.743990309495757
$ws.Range("L6").Value  = 52.401334422687093
$ws.Range("L7").Value  = 40.084286291781751
$ws.Range("L8").Value  = 58.6564425462321
$ws.Range("L9").Value  = 52.689880705632987
$ws.Range("L10").Value = 19.88866894869804
$ws.Range("L11").Value = 35.972443863264772
$ws.Range("L12").Value = 12.061786277026036
$ws.Range("L13").Value = -0.064288010286095529
$ws.Range("L14").Value = 34.132731805770057

# --- Formatting: column L should look like column K (same cell styles) ----
$ws.Range("K4:K14").Copy() | Out-Null
$ws.Range("L4:L14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row heights adjusted by Excel's auto-reflow after the edit -----------
$ws.Rows.Item(1).RowHeight = 67.5
for ($r = 4; $r -le 14; $r++) {
    $ws.Rows.Item($r).RowHeight = 14.25
}

# --- Clear the leftover "M7" selection from the source sheet view ---------
$ws.Range("A1").Select() | Out-Null
